$wb = $excel.ActiveWorkbook

# --- Sheet1 (DQ_Report): replace rows 2-11, clear old rows 12-17 ---
$ws1 = $wb.Worksheets.Item("DQ_Report")

# New data for rows 2-11 (A: PatientIdentifikator, B: ICD_Primärkode, C: Orpha_Kode, D: dq_msg)
$data = @(
    @("P_19285751", "E84.0",  587, "Relation  E84.0 - 587  ist im BfArM nicht vorhanden "),
    @("P_19285753", "E84.80", 587, "Relation  E84.80 - 587  ist im BfArM nicht vorhanden "),
    @("P_19285754", "E85.0",  586, "Relation  E85.0 - 586  ist im BfArM nicht vorhanden "),
    @("P_19285755", "E75.2",  325, "Relation  E75.2 - 325  ist im BfArM nicht vorhanden "),
    @("P_19285756", "E75.2",  320, "Relation  E75.2 - 320  ist im BfArM nicht vorhanden "),
    @("P_19285757", $null,    586, "Fehlendes ICD10 Code  "),
    @("P_19285758", $null,    587, "Orpha Kodierung  587  ist im BfArM-Mapping nicht enthalten Fehlendes ICD10 Code  "),
    @("P_19285759", "E75.2",  $null, "ICD10-Kodierung nicht eindeutig E75.2 Fehlendes Orpha_Kode  "),
    @("P_19285759", "E84.0",  $null, "Fehlendes Orpha_Kode  "),
    @("P_19285759", "D45",    $null, "Fehlendes Orpha_Kode  ")
)

# First clear out the old range (rows 2-17) so stale cells (e.g. old B/C values) don't linger
$ws1.Range("A2:D17").ClearContents()

$r = 2
foreach ($row in $data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    if ($null -eq $row[1]) {
        $ws1.Cells.Item($r, 2).Value = $null
    } else {
        $ws1.Cells.Item($r, 2).Value = $row[1]
    }
    if ($null -eq $row[2]) {
        $ws1.Cells.Item($r, 3).Value = $null
    } else {
        $ws1.Cells.Item($r, 3).Value = $row[2]
    }
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# --- Sheet2 (Statistik): update row 2 values ---
$ws2 = $wb.Worksheets.Item("Statistik")
$ws2.Cells.Item(2, 2).Value = 0.28
$ws2.Cells.Item(2, 3).Value = 99.72
$ws2.Cells.Item(2, 4).Value = 90.62
$ws2.Cells.Item(2, 5).Value = 98.3
$ws2.Cells.Item(2, 6).Value = 32
